$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (test cases 20160331 - 002 -> 003) appended below the
# existing log rows (rows 2-6), following the same column layout as row 1:
# A=Time, B=RunningTime(s), C=Preprocess, D=Features, E=Model,
# F=Model_Details, G=Test_Accuracy, H=Val_Accuracy, I=Template Filter,
# J=(unlabeled numeric column present in data rows)

$preprocess = 'trim "space" and ",", convert unicode to ascii, convert to lower, remove multiple spaces'
$features = '5 features: #ascii/(#ascii+#digit+#punctuation), #max_digit_skip_0_1, #(, #+, #/'
$model = 'Neuron Network'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300'
$templateFilter = '0 filters: '

# Parallel arrays (explicit ordering, no reliance on hashtable enumeration order)
$rowNums      = @(7, 8, 9, 10, 11)
$times        = @('20160401_015220', '20160401_015945', '20160401_020956', '20160401_021921', '20160401_022928')
$runningTimes = @(445.096, 610.582, 564.939, 607.277, 710.739)
$testAccs     = @(0.951333333333333, 0.957333333333333, 0.958, 0.958666666666667, 0.960666666666667)
$valAccs      = @(0.996699669966997, 0.996699669966997, 0.996699669966997, 0.996699669966997, 0.996699669966997)
$jVals        = @(0.357142857142857, 0.326530612244898, 0.326530612244898, 0.346938775510204, 0.336734693877551)

for ($i = 0; $i -lt $rowNums.Length; $i++) {
    $n = $rowNums[$i]
    $ws.Range("A$n").Value = $times[$i]
    $ws.Range("B$n").Value = $runningTimes[$i]
    $ws.Range("C$n").Value = $preprocess
    $ws.Range("D$n").Value = $features
    $ws.Range("E$n").Value = $model
    $ws.Range("F$n").Value = $modelDetails
    $ws.Range("G$n").Value = $testAccs[$i]
    $ws.Range("H$n").Value = $valAccs[$i]
    $ws.Range("I$n").Value = $templateFilter
    $ws.Range("J$n").Value = $jVals[$i]
}
